$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 42087
$ws.Range("F3").Value = 26
$ws.Range("F5").Value = 9730
$ws.Range("F6").Value = 208
$ws.Range("F7").Value = 935
$ws.Range("F8").Value = 935
$ws.Range("F9").Value = 742
$ws.Range("F10").Value = 224
$ws.Range("F11").Value = 309
$ws.Range("F12").Value = 953
$ws.Range("F14").Value = 763
$ws.Range("F15").Value = 326
$ws.Range("F16").Value = 1486
$ws.Range("F18").Value = 718
$ws.Range("F19").Value = 730
$ws.Range("F20").Value = 474
$ws.Range("F21").Value = 706
$ws.Range("F22").Value = 763
$ws.Range("F24").Value = 253
$ws.Range("F25").Value = 64
$ws.Range("F26").Value = 524
$ws.Range("F27").Value = 541
$ws.Range("F28").Value = 61
$ws.Range("F29").Value = 256
$ws.Range("F30").Value = 944
$ws.Range("F32").Value = 438
$ws.Range("F34").Value = 221
$ws.Range("F35").Value = 154
$ws.Range("F36").Value = 428
$ws.Range("F37").Value = 1321
$ws.Range("F38").Value = 305
$ws.Range("F39").Value = 1275
$ws.Range("F41").Value = 100
$ws.Range("F42").Value = 22
$ws.Range("F46").Value = 1
$ws.Range("F47").Value = 9

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 346
$ws.Range("F5").Value = 4458
$ws.Range("F7").Value = 340
$ws.Range("F11").Value = 135
$ws.Range("F19").Value = 6
$ws.Range("F20").Value = 4387

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2055
$ws.Range("F3").Value = 546
$ws.Range("F4").Value = 441

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2055
$ws.Range("F3").Value = 546
$ws.Range("F4").Value = 26
$ws.Range("F7").Value = 340
$ws.Range("F9").Value = 9730
$ws.Range("F10").Value = 208
$ws.Range("F11").Value = 935
$ws.Range("F12").Value = 935
$ws.Range("F14").Value = 441
$ws.Range("F15").Value = 935
$ws.Range("F16").Value = 135
$ws.Range("F17").Value = 224
$ws.Range("F18").Value = 309
$ws.Range("F19").Value = 953
$ws.Range("F22").Value = 763
$ws.Range("F23").Value = 326
$ws.Range("F24").Value = 1486
$ws.Range("F26").Value = 718
$ws.Range("F27").Value = 730
$ws.Range("F28").Value = 474
$ws.Range("F29").Value = 706
$ws.Range("F30").Value = 763
$ws.Range("F32").Value = 64
$ws.Range("F33").Value = 524
$ws.Range("F35").Value = 61
$ws.Range("F36").Value = 256
$ws.Range("F37").Value = 944
$ws.Range("F40").Value = 438
$ws.Range("F42").Value = 221
$ws.Range("F43").Value = 1275
$ws.Range("F45").Value = 100
$ws.Range("F49").Value = 9
